$d = $word.ActiveDocument

# Locate the "GearVR:" paragraph (currently the last paragraph in the body,
# holding the hidden _GoBack bookmark at its end).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "GearVR:*") {
        $target = $p
    }
}

# The _GoBack bookmark currently sits at the end of that paragraph; it needs
# to move to the end of the new final paragraph we are about to add, so
# remove it from its current spot first (it is a hidden bookmark, so it is
# not in the visible Bookmarks collection/count, but is still addressable
# by name).
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# Collapse a range to just before the paragraph mark of the GearVR
# paragraph so the inserted XML lands as new paragraphs right after it,
# without touching the existing run/text.
$r = $target.Range
$ins = $d.Range($r.End - 1, $r.End - 1)

$finalLabel = 'Final:'
$builtText = 'Built on Lab 02 to run on SteamVR with oculus touch controllers, refactored inputs for scalability, added the ambidexterity(Is that a word?) you can now grab the gun with either hand.'
$controlsText = 'Controls: X (Lower button on left controller) will reset the gun’s position to the counter, Z(lower button on right controller) will return you to the start screen, Left stick will move you (Note: It does not account for rotation of the headset, you move with forward always being towards the counter) Grip on either controller will attempt to pick up the gun with that hand, and the trigger on either will attempt to fire the gun if it is in that hand.'

$body = "<w:p/>" + `
        "<w:p><w:r><w:t>$finalLabel</w:t></w:r></w:p>" + `
        "<w:p><w:r><w:t>$builtText</w:t></w:r>" + `
        "<w:r><w:br/><w:t>$controlsText</w:t></w:r>" + `
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"

$xml = '<?xml version="1.0"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
       '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
       $body + `
       '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $ins.InsertXML($xml)
